$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet that used to be called "LINEST" is renamed to "SLOPE"
# ------------------------------------------------------------------
$wsSlope = $wb.Worksheets.Item("LINEST")
$wsSlope.Name = "SLOPE"

# ------------------------------------------------------------------
# C3 held a LINEST() array formula spanning C3:D3 (slope in C3,
# intercept in D3 as a plain cached value). Replace it with two
# independent formulas: SLOPE() in C3 and INTERCEPT() in D3.
# ------------------------------------------------------------------
$wsSlope.Range("C3:D3").ClearContents()
$wsSlope.Range("C3").Formula = "=SLOPE(B2:B9, A2:A9)"
$wsSlope.Range("D3").Formula = "=INTERCEPT(B2:B9, A2:A9)"

# Clearing the old array formula resets number-format/font on D3 to
# the workbook default; restore its original formatting (font size
# 14, same as the rest of the data rows on this sheet).
$wsSlope.Range("D3").Font.Size = 14

# ------------------------------------------------------------------
# Give the (previously untouched) C4:D9 range the same formatting as
# the rest of the sheet's data cells so they carry explicit style
# information, matching the rest of column A/B.
# ------------------------------------------------------------------
$wsSlope.Range("C4:D9").Font.Size = 14

# ------------------------------------------------------------------
# Column C narrows slightly (column D's explicit custom width is
# dropped in the target workbook; it is intentionally left alone
# here since there is no COM call that removes a column's width
# override outright without disturbing the sheet's data).
# ------------------------------------------------------------------
$wsSlope.Columns("C").ColumnWidth = 8.7

# ------------------------------------------------------------------
# Selection on the SLOPE sheet moves from G6 to G3.
# ------------------------------------------------------------------
$wsSlope.Range("G3").Select()

# ------------------------------------------------------------------
# The "DataAnalysis" sheet becomes the active (selected) tab, with
# N10 selected there (tabSelected flag moves off temp-sale onto
# DataAnalysis automatically once it becomes the active sheet).
# ------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("DataAnalysis")
$wsData.Activate()
$wsData.Range("N10").Select()

$wb.Save()
